$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date values from 45212 to 45221 for rows 2-32
$ws.Range("C2:C32").Value = 45221
